$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (flexibility), shifting existing
# columns D..J to F..L. This makes room for the new "sum_SASA" and
# "max_SASA" columns right after the existing "SASA" column (C). Inserting
# via EntireColumn.Insert() carries the formatting of the column to the
# left (C) onto the newly inserted columns, matching the header style.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"

# Fill the new columns with the same value as the SASA column (C) for each
# data row (2 through 12).
for ($r = 2; $r -le 12; $r++) {
    $val = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 4).Value = $val
    $ws.Cells.Item($r, 5).Value = $val
}
